$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3   = @{ Q = 58; R = 37 }
    10  = @{ Q = 31; R = 24 }
    17  = @{ Q = 12; R = 8 }
    23  = @{ Q = 21 }
    32  = @{ Q = 21; R = 15 }
    40  = @{ Q = 7;  R = 3 }
    49  = @{ Q = 16; R = 4 }
    58  = @{ Q = 84; R = 6 }
    66  = @{ Q = 42; R = 29 }
    74  = @{ Q = 38; R = 8 }
    78  = @{ Q = 56; R = 34 }
    89  = @{ Q = 2;  R = 1 }
    97  = @{ Q = 80; R = 48 }
    106 = @{ Q = 53; R = 10 }
    115 = @{ Q = 29; R = 9 }
    124 = @{ Q = 65; R = 25 }
    133 = @{ Q = 56; R = 44 }
    142 = @{ Q = 86; R = 35 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey('Q')) {
        $ws.Range("Q$row").Value = $vals['Q']
    }
    if ($vals.ContainsKey('R')) {
        $ws.Range("R$row").Value = $vals['R']
    }
}
